$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "96.443.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.703.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "651.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.42%  "

$ws.Range("E8").Value = "  -0.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

$ws.Range("E10").Value = "  -4.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.700.35"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000309"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +17.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.65%  "

$ws.Range("E15").Value = "  +2.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.391.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "96.298.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +13.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.680.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.504"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "520.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "

$ws.Range("E24").Value = "  -1.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000209"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.174"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.23%  "

$ws.Range("E30").Value = "  -1.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.186"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("E34").Value = "  +7.53%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -3.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "649.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.63%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.587"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.88"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.88%  "

$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.959"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.91%  "

$ws.Range("E46").Value = "  +1.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.429"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "

$ws.Range("E48").Value = "  -1.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("E51").Value = "  +2.11%  "
